# Append 9 new rows (22-30) to the "master-reg_center_user_h" sheet, following
# the same pattern as the existing rows: regcntr_id, usr_id, lang_code,
# is_active, cr_by, cr_dtimes, eff_dtimes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id (col A) and usr_id (col B) values for the new rows 22-30.
$regCntrIds = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$usrIds     = @(110021, 110022, 110023, 110024, 110025, 110026, 110027, 110028, 110029)

$startRow = 22
for ($i = 0; $i -lt $regCntrIds.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $regCntrIds[$i]
    $ws.Cells.Item($r, 2).Value = $usrIds[$i]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Select the first empty row below the data (full-row selection), mirroring
# the author's final on-screen selection state.
$ws.Range("A31:XFD1048576").Select()

# Configure the page setup for printing.
$ws.PageSetup.Orientation = 1
